# practicals/other_data_sets.xlsx — "Updated notes for day 2 and practical 4"
#
# Adds summary statistics (x bar bar, d2, Rbar) next to the "control 1"
# data table, and moves the active selection to C15 on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control 1")

# New labels (E7:E9) + values/formulas (F7:F9)
$ws.Range("E7").Value = "x bar bar"
$ws.Range("F7").Formula = "=AVERAGE(B2:B21)"

$ws.Range("E8").Value = "d2"
$ws.Range("F8").Value = 0.58

$ws.Range("E9").Value = "Rbar"
$ws.Range("F9").Formula = "=AVERAGE(C2:C21)"

# Match the author's final selection on this sheet
$ws.Activate()
$ws.Range("C15").Select() | Out-Null
